$d = $word.ActiveDocument

$d.Content.Find.Execute("28×83=2324", $true, $false, $false, $false, $false, $true, 1, $false, "98×69=6762", 2) | Out-Null
$d.Content.Find.Execute("96×39=3744", $true, $false, $false, $false, $false, $true, 1, $false, "80×83=6640", 2) | Out-Null
$d.Content.Find.Execute("62×36=2232", $true, $false, $false, $false, $false, $true, 1, $false, "40×94=3760", 2) | Out-Null
$d.Content.Find.Execute("22×32=704", $true, $false, $false, $false, $false, $true, 1, $false, "96×95=9120", 2) | Out-Null
$d.Content.Find.Execute("94×41=3854", $true, $false, $false, $false, $false, $true, 1, $false, "29×97=2813", 2) | Out-Null
$d.Content.Find.Execute("52×97=5044", $true, $false, $false, $false, $false, $true, 1, $false, "15×14=210", 2) | Out-Null
$d.Content.Find.Execute("39×84=3276", $true, $false, $false, $false, $false, $true, 1, $false, "33×16=528", 2) | Out-Null
$d.Content.Find.Execute("68×21=1428", $true, $false, $false, $false, $false, $true, 1, $false, "40×60=2400", 2) | Out-Null
$d.Content.Find.Execute("96×14=1344", $true, $false, $false, $false, $false, $true, 1, $false, "28×25=700", 2) | Out-Null
$d.Content.Find.Execute("53×45=2385", $true, $false, $false, $false, $false, $true, 1, $false, "11×88=968", 2) | Out-Null
$d.Content.Find.Execute("26×12=312", $true, $false, $false, $false, $false, $true, 1, $false, "80×23=1840", 2) | Out-Null
$d.Content.Find.Execute("49×84=4116", $true, $false, $false, $false, $false, $true, 1, $false, "21×27=567", 2) | Out-Null
$d.Content.Find.Execute("71×96=6816", $true, $false, $false, $false, $false, $true, 1, $false, "45×48=2160", 2) | Out-Null
$d.Content.Find.Execute("18×40=720", $true, $false, $false, $false, $false, $true, 1, $false, "39×15=585", 2) | Out-Null
$d.Content.Find.Execute("24×65=1560", $true, $false, $false, $false, $false, $true, 1, $false, "18×14=252", 2) | Out-Null
$d.Content.Find.Execute("17×37=629", $true, $false, $false, $false, $false, $true, 1, $false, "80×39=3120", 2) | Out-Null
$d.Content.Find.Execute("12×22=264", $true, $false, $false, $false, $false, $true, 1, $false, "73×66=4818", 2) | Out-Null
$d.Content.Find.Execute("69×84=5796", $true, $false, $false, $false, $false, $true, 1, $false, "73×63=4599", 2) | Out-Null
$d.Content.Find.Execute("58×60=3480", $true, $false, $false, $false, $false, $true, 1, $false, "92×20=1840", 2) | Out-Null
$d.Content.Find.Execute("38×99=3762", $true, $false, $false, $false, $false, $true, 1, $false, "36×98=3528", 2) | Out-Null
$d.Content.Find.Execute("39×34=1326", $true, $false, $false, $false, $false, $true, 1, $false, "82×52=4264", 2) | Out-Null
$d.Content.Find.Execute("29×71=2059", $true, $false, $false, $false, $false, $true, 1, $false, "68×23=1564", 2) | Out-Null
$d.Content.Find.Execute("12×54=648", $true, $false, $false, $false, $false, $true, 1, $false, "49×45=2205", 2) | Out-Null
$d.Content.Find.Execute("47×87=4089", $true, $false, $false, $false, $false, $true, 1, $false, "59×32=1888", 2) | Out-Null
$d.Content.Find.Execute("27×57=1539", $true, $false, $false, $false, $false, $true, 1, $false, "61×52=3172", 2) | Out-Null
